# Updated symbol list on Tue Dec 27 16:54:11 UTC 2022 with GitHub Actions
#
# Applies the scraped-price refresh described by the commit: a batch of
# "Price" (column D) updates, a couple of "Volume(1h)" label refreshes
# (column E) that pick up a Best/Worst-in-24h suffix, and a three-way
# re-ranking of rows 41-43 (BKEXToken / CEJI / KickToken) whose Coin name,
# Link and Price all rotate together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column is stored as *text* (t="inlineStr" in the source OOXML),
# not a number, even though the values look numeric. Excel's COM Value
# setter auto-converts a bare numeric-looking string into a real number, so
# we force text entry with a leading apostrophe and then reset the style
# back to "Normal" (the apostrophe trick stamps a quotePrefix style on the
# cell, which would otherwise show up as a spurious style change).
function Set-TextValue {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.Value = "'" + $Text
    $cell.Style = "Normal"
}

# ---- Column D: Price refresh ----
Set-TextValue "D2"  "245.53"
Set-TextValue "D3"  "23.71"
Set-TextValue "D4"  "5.369"
Set-TextValue "D5"  "0.05878"
Set-TextValue "D6"  "3.376"
Set-TextValue "D7"  "6.480"
Set-TextValue "D9"  "0.9207"
Set-TextValue "D11" "0.07416"
Set-TextValue "D13" "0.03071"
Set-TextValue "D14" "0.09362"
Set-TextValue "D15" "3.858"
Set-TextValue "D16" "0.001550"
Set-TextValue "D17" "0.04721"
Set-TextValue "D18" "0.0005992"
Set-TextValue "D19" "0.005950"
Set-TextValue "D20" "0.001245"
Set-TextValue "D21" "0.004715"
Set-TextValue "D27" "0.0002654"
Set-TextValue "D40" "0.03869"
Set-TextValue "D44" "0.008467"
Set-TextValue "D45" "0.00005251"
Set-TextValue "D47" "0.7103"
Set-TextValue "D48" "0.001738"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.0002001"

# ---- Column E: Volume(1h) label gains a Best/Worst-in-24h suffix ----
$ws.Range("E18").Value = "17OneONEWorstin24h"

# ---- Rows 41-43: BKEXToken / CEJI / KickToken rotate rank positions ----
# Row 41 (was BKEXToken) becomes KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006343"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42 (was CEJI) becomes BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1066"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43 (was KickToken) becomes CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.003201"
$ws.Range("E43").Value = "42CEJICEJI"
